# Backlog update: "Se agrego un objeto game que sirve como facade del juego."
#
# - Hoja1: la tarea "Start"/"Despacahar las urls" se marca como terminada (ok),
#   y se agrega una nueva tarea: "Notificar a los jugadores" / "A partir del
#   player manager notificar a cada jugador con lo que corresponda".
# - "Agregar jugadores": se corrige el texto "generar clase que permisista los
#   jugadores" -> "generar clase que persista los jugadores".
# - La hoja activa pasa de "Agregar jugadores" a "Hoja1", con la selección
#   parada justo debajo de la fila nueva.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Agregar jugadores")

# Hoja1: la tarea de la fila 13 ("Despacahar las urls") queda lista.
$ws1.Range("D13").Value = "ok"

# Hoja1: nueva fila 14 con la siguiente tarea del backlog.
$ws1.Range("B14").Value = "Notificar a los jugadores"
$ws1.Range("C14").Value = "A partir del player manager notificar a cada jugador con lo que corresponda"

# "Agregar jugadores": arreglar el typo "permisista" -> "persista".
$ws2.Range("B5").Value = "generar clase que persista los jugadores"

# Selección / hoja activa: al terminar, el foco vuelve a Hoja1.
$null = $ws2.Range("B4").Select()
$null = $ws1.Activate()
$null = $ws1.Range("C15").Select()
